$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 930.1111
$ws.Range("I6").Value = 228.5
$ws.Range("J6").Value = 2333.3333
$ws.Range("K6").Value = 685.5
$ws.Range("L6").Value = 6999.999899999999
$ws.Range("M6").Value = -573.5
$ws.Range("N6").Value = -7223.999899999999

$ws.Range("H8").Value = 668.4286
$ws.Range("I8").Value = 32.545456
$ws.Range("K8").Value = 97.636368
$ws.Range("M8").Value = 41.363632

$ws.Range("H9").Value = 376.69565
$ws.Range("I9").Value = 137.27272
$ws.Range("J9").Value = 596.1667
$ws.Range("K9").Value = 137.27272
$ws.Range("L9").Value = 596.1667
$ws.Range("M9").Value = 31.72728000000001
$ws.Range("N9").Value = -934.1667

$ws.Range("H12").Value = 285914.28
$ws.Range("I12").Value = 116.666664
$ws.Range("J12").Value = 500262.5
$ws.Range("K12").Value = 116.666664
$ws.Range("L12").Value = 500262.5
$ws.Range("M12").Value = 53.333336
$ws.Range("N12").Value = -500602.5

$ws.Range("H21").Value = 16669.834
$ws.Range("I21").Value = 18003.8
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 18003.8
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -17535.8
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 16669.834
$ws.Range("I23").Value = 18003.8
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 18003.8
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -17769.8
$ws.Range("N23").Value = -10468

$ws.Range("H33").Value = 614.7917
$ws.Range("I33").Value = 456.55554
$ws.Range("J33").Value = 1089.5
$ws.Range("K33").Value = 456.55554
$ws.Range("L33").Value = 1089.5
$ws.Range("M33").Value = -227.55554
$ws.Range("N33").Value = -1547.5

$ws.Range("H38").Value = 853.94446
$ws.Range("I38").Value = 216.54546
$ws.Range("J38").Value = 1855.5714
$ws.Range("K38").Value = 649.6363799999999
$ws.Range("L38").Value = 5566.7142
$ws.Range("M38").Value = -277.6363799999999
$ws.Range("N38").Value = -6310.7142

$ws.Range("H99").Value = 709.7778
$ws.Range("I99").Value = 369.7143
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 1109.1429
$ws.Range("L99").Value = 5700
$ws.Range("M99").Value = 388.8571000000002
$ws.Range("N99").Value = -8696

$ws.Range("H112").Value = 4167820.2
$ws.Range("I112").Value = 2055.7144
$ws.Range("J112").Value = 4718015.5
$ws.Range("K112").Value = 6167.1432
$ws.Range("L112").Value = 14154046.5
$ws.Range("M112").Value = -5059.1432
$ws.Range("N112").Value = -14156262.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1584.1666
$ws.Range("I45").Value = 1019.25
$ws.Range("J45").Value = 6103.5
$ws.Range("K45").Value = 1019.25
$ws.Range("L45").Value = 6103.5
$ws.Range("M45").Value = -642.25
$ws.Range("N45").Value = -6857.5

$ws.Range("H135").Value = 39000
$ws.Range("J135").Value = 39000
$ws.Range("L135").Value = 39000
$ws.Range("N135").Value = -49140

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3242.6
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920

$ws.Range("H134").Value = 2320.456
$ws.Range("I134").Value = 2080.3555
$ws.Range("J134").Value = 3220.8333
$ws.Range("K134").Value = 6241.066500000001
$ws.Range("L134").Value = 9662.499899999999
$ws.Range("M134").Value = -3706.066500000001
$ws.Range("N134").Value = -14732.4999

$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

$ws.Range("H137").Value = 29362.23
$ws.Range("J137").Value = 29166.666
$ws.Range("L137").Value = 29166.666
$ws.Range("N137").Value = -39366.666

$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13891368
$ws.Range("I58").Value = 1431.125
$ws.Range("J58").Value = 41671240
$ws.Range("K58").Value = 1431.125
$ws.Range("L58").Value = 41671240
$ws.Range("M58").Value = -1228.125
$ws.Range("N58").Value = -41671646

$ws.Range("H105").Value = 3578.3076
$ws.Range("I105").Value = 3350.8
$ws.Range("J105").Value = 4336.6665
$ws.Range("K105").Value = 3350.8
$ws.Range("L105").Value = 4336.6665
$ws.Range("M105").Value = -1603.8
$ws.Range("N105").Value = -7830.6665

$ws.Range("H132").Value = 2812.4167
$ws.Range("I132").Value = 1625.4
$ws.Range("J132").Value = 5510.1816
$ws.Range("K132").Value = 4876.200000000001
$ws.Range("L132").Value = 16530.5448
$ws.Range("M132").Value = -2346.200000000001
$ws.Range("N132").Value = -21590.5448

$ws.Range("H134").Value = 1650.9032
$ws.Range("I134").Value = 1282.7115
$ws.Range("K134").Value = 3848.1345
$ws.Range("M134").Value = -1313.1345

$ws.Range("H136").Value = 13891368
$ws.Range("I136").Value = 1431.125
$ws.Range("J136").Value = 41671240
$ws.Range("K136").Value = 4293.375
$ws.Range("L136").Value = 125013720
$ws.Range("M136").Value = -1743.375
$ws.Range("N136").Value = -125018820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1867.0588
$ws.Range("I131").Value = 1908.75
$ws.Range("J131").Value = 1854.2307
$ws.Range("K131").Value = 5726.25
$ws.Range("L131").Value = 5562.6921
$ws.Range("M131").Value = -686.25
$ws.Range("N131").Value = -15642.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 988.5333000000001
$ws.Range("I107").Value = 247.11111
$ws.Range("J107").Value = 2100.6667
$ws.Range("K107").Value = 247.11111
$ws.Range("L107").Value = 2100.6667
$ws.Range("M107").Value = 1672.88889
$ws.Range("N107").Value = -5940.6667

$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2710.8333
$ws.Range("I122").Value = 2482.7585
$ws.Range("J122").Value = 3655.7144
$ws.Range("K122").Value = 7448.2755
$ws.Range("L122").Value = 10967.1432
$ws.Range("M122").Value = -4998.2755
$ws.Range("N122").Value = -15867.1432

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 925
$ws.Range("J81").Value = 1450
$ws.Range("L81").Value = 2900
$ws.Range("N81").Value = -5022

$ws.Range("H84").Value = 925
$ws.Range("J84").Value = 1450
$ws.Range("L84").Value = 14500
$ws.Range("N84").Value = -25108

$ws.Range("H122").Value = 258171.92
$ws.Range("I122").Value = 313618.9
$ws.Range("K122").Value = 940856.7000000001
$ws.Range("M122").Value = -938406.7000000001

$ws.Range("H136").Value = 1588.4
$ws.Range("I136").Value = 738.8929000000001
$ws.Range("J136").Value = 4986.4287
$ws.Range("K136").Value = 2216.6787
$ws.Range("L136").Value = 14959.2861
$ws.Range("M136").Value = 333.3212999999996
$ws.Range("N136").Value = -20059.2861

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
